# Apply "Block IDs in CAN Library" edit:
# 1. Update the "main" sheet: rename the G-column header to "Rx Any (T/F)",
#    populate the new Rx-Any flag column, and clear out the old
#    Set-Address / IVT current-voltage rows (that data moves to a new sheet).
# 2. Add a new "identifiers" worksheet listing the IVT CAN identifiers.
# 3. Restore view state (zoom / selection) to match the edited workbook.

$wb = $excel.ActiveWorkbook
$main = $wb.Worksheets.Item("main")

# --- Header: column G is now "Rx Any (T/F)" instead of "Set Address" ---
$main.Range("G1").Value2 = "Rx Any (T/F)"

# --- Populate the new Rx Any (T/F) column ---
$main.Range("G2").Value2 = "F"
$main.Range("G4").Value2 = "T"
$main.Range("G6").Value2 = "F"
$main.Range("G7").Value2 = "F"
$main.Range("G11").Value2 = "F"

# --- Remove the old Set Address / Current / Voltage 1-3 rows; that
#     information now lives on the "identifiers" sheet ---
$main.Range("C12").ClearContents() | Out-Null
$main.Range("E12").ClearContents() | Out-Null
$main.Range("F12").ClearContents() | Out-Null
$main.Range("E13").ClearContents() | Out-Null
$main.Range("F13").ClearContents() | Out-Null
$main.Range("E14").ClearContents() | Out-Null
$main.Range("F14").ClearContents() | Out-Null
$main.Range("E15").ClearContents() | Out-Null
$main.Range("F15").ClearContents() | Out-Null

# --- Add the new "identifiers" sheet after "main" ---
$idSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $main)
$idSheet.Name = "identifiers"

$idSheet.Range("A1").Value2 = "Name"
$idSheet.Range("B1").Value2 = "ID"
$idSheet.Range("C1").Value2 = "Data Type"
$idSheet.Range("D1").Value2 = "Receive (T/F)"

$idSheet.Range("A2").Value2 = "IVT_Current"
$idSheet.Range("B2").Value2 = "0x521"
$idSheet.Range("C2").Value2 = "uint8_t[6]"
$idSheet.Range("D2").Value2 = "T"

$idSheet.Range("A3").Value2 = "IVT_Voltage1"
$idSheet.Range("B3").Value2 = "0x522"
$idSheet.Range("C3").Value2 = "uint8_t[6]"
$idSheet.Range("D3").Value2 = "T"

$idSheet.Range("A4").Value2 = "IVT_Voltage2"
$idSheet.Range("B4").Value2 = "0x523"
$idSheet.Range("C4").Value2 = "uint8_t[6]"
$idSheet.Range("D4").Value2 = "T"

$idSheet.Range("A5").Value2 = "IVT_Voltage3"
$idSheet.Range("B5").Value2 = "0x524"
$idSheet.Range("C5").Value2 = "uint8_t[6]"
$idSheet.Range("D5").Value2 = "T"

$idSheet.Columns.Item(1).AutoFit() | Out-Null
$idSheet.Columns.Item(4).AutoFit() | Out-Null

$idSheet.Range("B3").Select() | Out-Null

# --- Restore "main" as the active sheet/view ---
$main.Activate()
$excel.ActiveWindow.Zoom = 125
$main.Range("G11").Select() | Out-Null
